$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-of dates) ---
$ws.Cells.Item(8, 1).Value = "Volume 30   Number  15"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  4/10/2023  Through  4/16/2023"

# --- Cells that flip from the text "n/a" placeholder to a real number ---
# (copy number-style formatting from a donor cell first, then set the value)
$ws.Cells.Item(15, 6).Copy($ws.Cells.Item(15, 4))
$ws.Cells.Item(15, 4).Value = 2
$ws.Cells.Item(16, 5).Copy($ws.Cells.Item(15, 5))
$ws.Cells.Item(15, 5).Value = -100
$ws.Cells.Item(15, 6).Copy($ws.Cells.Item(15, 7))
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(16, 5).Copy($ws.Cells.Item(15, 8))
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(20, 6).Copy($ws.Cells.Item(20, 4))
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(16, 5).Copy($ws.Cells.Item(20, 5))
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Copy($ws.Cells.Item(26, 3))
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(20, 6).Copy($ws.Cells.Item(27, 4))
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(16, 5).Copy($ws.Cells.Item(27, 5))
$ws.Cells.Item(27, 5).Value = 300

# --- Cells that flip from a real number to the text "n/a" placeholder ---
# (copy full cell - format & shared-string value - from a donor cell that already holds it)
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(14, 7))
$ws.Cells.Item(14, 5).Copy($ws.Cells.Item(14, 8))
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(22, 7))
$ws.Cells.Item(14, 5).Copy($ws.Cells.Item(22, 8))
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(30, 7))
$ws.Cells.Item(14, 5).Copy($ws.Cells.Item(30, 8))

# --- Plain numeric value updates ---
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(15, 9).Value = 7
$ws.Cells.Item(15, 10).Value = 6
$ws.Cells.Item(15, 11).Value = 16.666666666666
$ws.Cells.Item(15, 12).Value = -12.5
$ws.Cells.Item(15, 13).Value = 250
$ws.Cells.Item(15, 14).Value = -50
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 400
$ws.Cells.Item(16, 6).Value = 17
$ws.Cells.Item(16, 7).Value = 9
$ws.Cells.Item(16, 8).Value = 88.888888888888
$ws.Cells.Item(16, 9).Value = 48
$ws.Cells.Item(16, 10).Value = 78
$ws.Cells.Item(16, 11).Value = -38.461538461538
$ws.Cells.Item(16, 12).Value = 29.729729729729
$ws.Cells.Item(16, 13).Value = 23.076923076923
$ws.Cells.Item(16, 14).Value = -81.322957198443
$ws.Cells.Item(17, 3).Value = 5
$ws.Cells.Item(17, 4).Value = 6
$ws.Cells.Item(17, 5).Value = -16.666666666666
$ws.Cells.Item(17, 6).Value = 23
$ws.Cells.Item(17, 7).Value = 17
$ws.Cells.Item(17, 8).Value = 35.294117647058
$ws.Cells.Item(17, 9).Value = 64
$ws.Cells.Item(17, 10).Value = 56
$ws.Cells.Item(17, 11).Value = 14.285714285714
$ws.Cells.Item(17, 12).Value = 52.380952380952
$ws.Cells.Item(17, 13).Value = 48.837209302325
$ws.Cells.Item(17, 14).Value = -56.462585034013
$ws.Cells.Item(18, 3).Value = 6
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(18, 5).Value = 100
$ws.Cells.Item(18, 7).Value = 25
$ws.Cells.Item(18, 8).Value = -20
$ws.Cells.Item(18, 9).Value = 80
$ws.Cells.Item(18, 10).Value = 101
$ws.Cells.Item(18, 11).Value = -20.79207920792
$ws.Cells.Item(18, 12).Value = -10.112359550561
$ws.Cells.Item(18, 13).Value = 6.666666666666
$ws.Cells.Item(18, 14).Value = -67.611336032388
$ws.Cells.Item(19, 3).Value = 18
$ws.Cells.Item(19, 4).Value = 20
$ws.Cells.Item(19, 5).Value = -10
$ws.Cells.Item(19, 6).Value = 81
$ws.Cells.Item(19, 7).Value = 75
$ws.Cells.Item(19, 8).Value = 8
$ws.Cells.Item(19, 9).Value = 274
$ws.Cells.Item(19, 10).Value = 274
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 80.263157894736
$ws.Cells.Item(19, 13).Value = 20.704845814978
$ws.Cells.Item(19, 14).Value = -33.007334963325
$ws.Cells.Item(20, 6).Value = 5
$ws.Cells.Item(20, 8).Value = 150
$ws.Cells.Item(20, 9).Value = 11
$ws.Cells.Item(20, 10).Value = 12
$ws.Cells.Item(20, 11).Value = -8.333333333333
$ws.Cells.Item(20, 12).Value = -26.666666666666
$ws.Cells.Item(20, 13).Value = -26.666666666666
$ws.Cells.Item(20, 14).Value = -92.413793103448
$ws.Cells.Item(21, 3).Value = 35
$ws.Cells.Item(21, 4).Value = 33
$ws.Cells.Item(21, 5).Value = 6.060606060606
$ws.Cells.Item(21, 6).Value = 148
$ws.Cells.Item(21, 7).Value = 130
$ws.Cells.Item(21, 8).Value = 13.846153846153
$ws.Cells.Item(21, 9).Value = 484
$ws.Cells.Item(21, 10).Value = 529
$ws.Cells.Item(21, 11).Value = -8.506616257088
$ws.Cells.Item(21, 12).Value = 40.289855072463
$ws.Cells.Item(21, 13).Value = 20.398009950248
$ws.Cells.Item(21, 14).Value = -60.425183973834
$ws.Cells.Item(23, 3).Value = 4
$ws.Cells.Item(23, 4).Value = 4
$ws.Cells.Item(23, 6).Value = 15
$ws.Cells.Item(23, 7).Value = 12
$ws.Cells.Item(23, 8).Value = 25
$ws.Cells.Item(23, 9).Value = 38
$ws.Cells.Item(23, 10).Value = 49
$ws.Cells.Item(23, 11).Value = -22.448979591836
$ws.Cells.Item(23, 12).Value = -34.482758620689
$ws.Cells.Item(23, 13).Value = 18.75
$ws.Cells.Item(24, 3).Value = 27
$ws.Cells.Item(24, 4).Value = 44
$ws.Cells.Item(24, 5).Value = -38.636363636363
$ws.Cells.Item(24, 6).Value = 102
$ws.Cells.Item(24, 7).Value = 200
$ws.Cells.Item(24, 8).Value = -49
$ws.Cells.Item(24, 9).Value = 413
$ws.Cells.Item(24, 10).Value = 615
$ws.Cells.Item(24, 11).Value = -32.845528455284
$ws.Cells.Item(24, 12).Value = 75
$ws.Cells.Item(24, 13).Value = -8.830022075055
$ws.Cells.Item(25, 3).Value = 9
$ws.Cells.Item(25, 4).Value = 11
$ws.Cells.Item(25, 5).Value = -18.181818181818
$ws.Cells.Item(25, 6).Value = 49
$ws.Cells.Item(25, 7).Value = 34
$ws.Cells.Item(25, 8).Value = 44.117647058823
$ws.Cells.Item(25, 9).Value = 134
$ws.Cells.Item(25, 10).Value = 129
$ws.Cells.Item(25, 11).Value = 3.875968992248
$ws.Cells.Item(25, 12).Value = 63.414634146341
$ws.Cells.Item(25, 13).Value = 2.290076335877
$ws.Cells.Item(26, 4).Value = 4
$ws.Cells.Item(26, 5).Value = -75
$ws.Cells.Item(26, 6).Value = 3
$ws.Cells.Item(26, 7).Value = 5
$ws.Cells.Item(26, 8).Value = -40
$ws.Cells.Item(26, 9).Value = 13
$ws.Cells.Item(26, 10).Value = 13
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 8.333333333333
$ws.Cells.Item(27, 3).Value = 4
$ws.Cells.Item(27, 6).Value = 6
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 20
$ws.Cells.Item(27, 10).Value = 26
$ws.Cells.Item(27, 11).Value = -23.076923076923
$ws.Cells.Item(27, 12).Value = 53.846153846153
$ws.Cells.Item(28, 14).Value = -91.666666666666
$ws.Cells.Item(29, 14).Value = -87.5
